$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 30
$ws.Range("C2").Value = 45
$ws.Range("I2").Value = 10

$ws.Range("B3").Value = 60
$ws.Range("C3").Value = 28
$ws.Range("G3").Value = 100
$ws.Range("I3").Value = 20

$ws.Range("D4").Value = 90
$ws.Range("E4").Value = 150
$ws.Range("G4").Value = 37
$ws.Range("I4").Value = 130

$ws.Range("I5").Select

$excel.ActiveWindow.Width = 15345
$excel.ActiveWindow.Height = 4470
